$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header labels (row 1) - keep bold styling
$ws.Range("A1").Value2 = "browser"
$ws.Range("B1").Value2 = "url"

# Swap the data row (row 2)
$ws.Range("A2").Value2 = "edge"
$ws.Range("B2").Value2 = "http://www.way2automation.com/angularjs-protractor/webtables/"

# Turn B2 into a hyperlink pointing at itself, applying the built-in Hyperlink style
$ws.Hyperlinks.Add($ws.Range("B2"), "http://www.way2automation.com/angularjs-protractor/webtables/") | Out-Null

# Update the active selection to match the authored state
$ws.Range("D7").Select() | Out-Null
